$d = $word.ActiveDocument

# Locate the paragraph that starts with the MIPR/7600B label
# (Paragraph 2: "Military Interdepartmental Purchase Request (MIPR) or 7600B(s): {fundingDocInfo}")
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Military Interdepartmental Purchase Request*") {
        $targetPara = $para
        break
    }
}

if ($targetPara -ne $null) {
    $searchRange = $targetPara.Range

    # Remove the leading label text "Military Interdepartmental Purchase Request (MIPR) or 7600B(s): "
    # so the paragraph starts directly with the {fundingDocInfo} placeholder.
    $searchRange.Find.ClearFormatting()
    $searchRange.Find.Execute(
        "Military Interdepartmental Purchase Request (MIPR) or 7600B(s): ",
        $true, $false, $false, $false, $false, $true, 1, $false, "", 2
    ) | Out-Null

    # Drop the paragraph's small left indent (w:ind w:left="27") now that the
    # leading label text is gone.
    $targetPara.Range.ParagraphFormat.LeftIndent = 0
}
